$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# New handoff report was generated for the f920e240-... file (row 7 on every
# sheet): the handoff timestamps advance to the latest run.
$wsOverview.Range("D7").Value = "2016-03-22 22:46:56"
$wsZhCn.Range("E7").Value = "2016-03-22 22:46:52"
$wsDeDe.Range("E7").Value = "2016-03-22 22:46:56"
